$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest cryptos snapshot.
# D-column values are forced to Text (leading apostrophe) to match the
# source data which stores prices like "24.700.50" / "316.65" as plain text.

$ws.Range("D2").Value = "'24.700.50"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "'1.694.37"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").Value = "'316.65"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'0.3955"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'0.4070"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").Value = "'1.493"
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("D10").Value = "'1.003"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "'52.92"
$ws.Range("E11").Value = "  -6.63%  "
$ws.Range("D12").Value = "'0.08942"
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("D13").Value = "'7.277"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").Value = "'23.54"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "'8.031"
$ws.Range("E15").Value = "  +5.65%  "
$ws.Range("D16").Value = "'0.00001327"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "'1.698.20"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "'100.00"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "'0.07033"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "'19.67"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").Value = "'6.989"
$ws.Range("E21").Value = "  +3.91%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "'14.33"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").Value = "'24.676.12"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'3.270"
$ws.Range("E25").Value = "  +8.03%  "
$ws.Range("D26").Value = "'2.364"
$ws.Range("E26").Value = "  +2.34%  "
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("D28").Value = "'162.22"
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("D29").Value = "'136.26"
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("D30").Value = "'5.160"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("D31").Value = "'7.503"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").Value = "'0.08689"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("D33").Value = "'1.053"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("D34").Value = "'7.074"
$ws.Range("E34").Value = "  -4.40%  "
$ws.Range("D35").Value = "'11.43"
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("D36").Value = "'0.2741"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").Value = "'1.887"
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("D38").Value = "'14.46"
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("D39").Value = "'0.09248"
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("D40").Value = "'0.02727"
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("D41").Value = "'1.473"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").Value = "'0.7667"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "'16.16"
$ws.Range("E43").Value = "  +5.11%  "
$ws.Range("D44").Value = "'2.591"
$ws.Range("E44").Value = "  +5.24%  "
$ws.Range("D45").Value = "'0.7180"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("D48").Value = "'140.26"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "'1.322"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").Value = "'91.03"
$ws.Range("E50").Value = "  +5.31%  "
$ws.Range("D51").Value = "'0.07981"
$ws.Range("E51").Value = "  -0.23%  "

# Re-apply the default "Normal" style to the Price column so the quote-prefix
# (text-forcing) marker does not leave a stray explicit cell style behind -
# matches the unstyled inlineStr cells in the original workbook.
$ws.Range("D2:D51").Style = "Normal"
